$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo" moves from D to E)
$ws.Columns.Item(4).Insert()

# New header for column D
$ws.Range("D1").Value = "MAE"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "MAE"

# New value for column D, row 2
$ws.Range("D2").Value = 2.046624516121589
